# Finalise the SDMX models for the four new datasets.
# This script updates the DSD concept table with renamed / new concepts,
# updates the Indicator label sheet, and fixes up the active-sheet/selection
# state to match the final authored workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. DSD sheet: rename / extend the concept scheme rows
# ---------------------------------------------------------------------
$dsd = $wb.Worksheets.Item("DSD")

# REF_AREA -> GEO_PICT (id renamed, everything else about the row unchanged)
$dsd.Range("A4").Value = "GEO_PICT"

# UNIT_MEASURE now points at the common codelist
$dsd.Range("F12").Value = "CL_COM_UNIT_MEASURE"

# New row: UNIT_MULT (inserted where OBS_STATUS used to be; rows below shift down)
$dsd.Range("A13").Value = "UNIT_MULT"
$dsd.Range("B13").Value = "Unit multiplier"
$dsd.Range("C13").Value = "Y"
$dsd.Range("D13").Value = "Attribute"
$dsd.Range("E13").Value = "Coded"
$dsd.Range("F13").Value = "CL_COM_UNIT_MULT"
$dsd.Range("G13").Value = "Y"

# OBS_STATUS moves down to row 14 and now points at the common codelist
$dsd.Range("A14").Value = "OBS_STATUS"
$dsd.Range("B14").Value = "Observation Status"
$dsd.Range("C14").Value = "Y"
$dsd.Range("D14").Value = "Attribute"
$dsd.Range("E14").Value = "Coded"
$dsd.Range("F14").Value = "CL_COM_OBS_STATUS"
$dsd.Range("G14").Value = "Y"

# New row: DATA_SOURCE
$dsd.Range("A15").Value = "DATA_SOURCE"
$dsd.Range("B15").Value = "Data source"
$dsd.Range("C15").Value = "N"
$dsd.Range("D15").Value = "Attribute"
$dsd.Range("E15").Value = "Uncoded"
$dsd.Range("F15").Value = "Text"
$dsd.Range("G15").Value = "N"

# COMMENT -> OBS_COMMENT (id renamed), moved down to row 16
$dsd.Range("A16").Value = "OBS_COMMENT"
$dsd.Range("B16").Value = "Comments"
$dsd.Range("C16").Value = "Y"
$dsd.Range("D16").Value = "Attribute"
$dsd.Range("E16").Value = "Uncoded"
$dsd.Range("F16").Value = "Text"
$dsd.Range("G16").Value = "Y"

# New row: CONF_STATUS
$dsd.Range("A17").Value = "CONF_STATUS"
$dsd.Range("B17").Value = "Confidentiality status"
$dsd.Range("C17").Value = "Y"
$dsd.Range("D17").Value = "Attribute"
$dsd.Range("E17").Value = "Coded"
$dsd.Range("F17").Value = "CL_COM_CONF_STATUS"
$dsd.Range("G17").Value = "Y"

# ---------------------------------------------------------------------
# 2. Indicator sheet: the "N" id row now reuses the label text as its id
# ---------------------------------------------------------------------
$indicator = $wb.Worksheets.Item("Indicator")
$indicator.Range("A2").Value = "Number of households"

# ---------------------------------------------------------------------
# 3. Fix up active sheet / selection so DSD (now finished) is the one
#    shown & selected when the workbook is opened, matching the final
#    authoring session.
# ---------------------------------------------------------------------
$tuber = $wb.Worksheets.Item("AGRICULTURE_TUBER")
$tuber.Range("F10").Select() | Out-Null

$dsd.Activate()
$dsd.Range("A12:G17").Select() | Out-Null
